$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = -0.1768982652262758
$ws.Cells.Item(2, 3).Value = 0.6795443542960273
$ws.Cells.Item(2, 4).Value = 1.127648372139759
$ws.Cells.Item(2, 5).Value = 1.061907892493393
$ws.Cells.Item(2, 6).Value = 1.071710210249048
$ws.Cells.Item(2, 7).Value = 22

$ws.Cells.Item(3, 2).Value = 0.6385729239511911
$ws.Cells.Item(3, 3).Value = 0.9114413858617648
$ws.Cells.Item(3, 4).Value = 2.040117013382039
$ws.Cells.Item(3, 5).Value = 1.428326647998293
$ws.Cells.Item(3, 6).Value = 1.309182460884421
$ws.Cells.Item(3, 7).Value = 21

$ws.Cells.Item(4, 2).Value = 0.535249571191142
$ws.Cells.Item(4, 3).Value = 1.202902938844514
$ws.Cells.Item(4, 4).Value = 3.603967812443823
$ws.Cells.Item(4, 5).Value = 1.898411918537129
$ws.Cells.Item(4, 6).Value = 1.868710703577968
$ws.Cells.Item(4, 7).Value = 20

$ws.Cells.Item(5, 2).Value = 0.6699125474792575
$ws.Cells.Item(5, 3).Value = 0.8506259737546845
$ws.Cells.Item(5, 4).Value = 1.144803732585156
$ws.Cells.Item(5, 5).Value = 1.069955014281047
$ws.Cells.Item(5, 6).Value = 0.8571398600703368
$ws.Cells.Item(5, 7).Value = 19

$ws.Cells.Item(6, 2).Value = 0.5261311751648745
$ws.Cells.Item(6, 3).Value = 0.6862967174377638
$ws.Cells.Item(6, 4).Value = 0.7267385007115407
$ws.Cells.Item(6, 5).Value = 0.8524895897965796
$ws.Cells.Item(6, 6).Value = 0.6902105718828747
$ws.Cells.Item(6, 7).Value = 18

$ws.Cells.Item(7, 2).Value = 0.3302666592216413
$ws.Cells.Item(7, 3).Value = 0.5314189662411473
$ws.Cells.Item(7, 4).Value = 0.3989061542592715
$ws.Cells.Item(7, 5).Value = 0.631590179039598
$ws.Cells.Item(7, 6).Value = 0.5549274444195054
$ws.Cells.Item(7, 7).Value = 17

$ws.Cells.Item(8, 2).Value = 0.3524319253953514
$ws.Cells.Item(8, 3).Value = 0.4057300614388591
$ws.Cells.Item(8, 4).Value = 0.2188174786034947
$ws.Cells.Item(8, 5).Value = 0.4677793054459493
$ws.Cells.Item(8, 6).Value = 0.3176735709340561
$ws.Cells.Item(8, 7).Value = 16
